$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 577
$ws.Range("I6").Value = 96.25
$ws.Range("K6").Value = 288.75
$ws.Range("M6").Value = -176.75
$ws.Range("H10").Value = 33202
$ws.Range("J10").Value = 33202
$ws.Range("L10").Value = 33202
$ws.Range("N10").Value = -33788
$ws.Range("H132").Value = 3863370.5
$ws.Range("I132").Value = 4466522
$ws.Range("J132").Value = 3201.2
$ws.Range("K132").Value = 13399566
$ws.Range("L132").Value = 9603.599999999999
$ws.Range("M132").Value = -13397036
$ws.Range("N132").Value = -14663.6
$ws.Range("H137").Value = 1460.6923
$ws.Range("I137").Value = 1236.4286
$ws.Range("K137").Value = 3709.2858
$ws.Range("M137").Value = -1159.2858
$ws.Range("H138").Value = 1630.12
$ws.Range("I138").Value = 1199.619
$ws.Range("J138").Value = 2178.0303
$ws.Range("K138").Value = 3598.857
$ws.Range("L138").Value = 6534.090899999999
$ws.Range("M138").Value = 1541.143
$ws.Range("N138").Value = -16814.0909

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I10").Value = 4500
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 4500
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -4330
$ws.Range("N10").ClearContents()
$ws.Range("H32").Value = 14576.553
$ws.Range("I32").Value = 14967.012
$ws.Range("J32").Value = 10888.889
$ws.Range("K32").Value = 14967.012
$ws.Range("L32").Value = 10888.889
$ws.Range("M32").Value = -14680.012
$ws.Range("N32").Value = -11462.889
$ws.Range("H61").Value = 1356.3959
$ws.Range("I61").Value = 1088.0322
$ws.Range("J61").Value = 1845.7646
$ws.Range("K61").Value = 1088.0322
$ws.Range("L61").Value = 1845.7646
$ws.Range("M61").Value = -876.0322000000001
$ws.Range("N61").Value = -2269.7646
$ws.Range("H74").Value = 939.1
$ws.Range("I74").Value = 807.1429000000001
$ws.Range("J74").Value = 1862.8
$ws.Range("K74").Value = 807.1429000000001
$ws.Range("L74").Value = 1862.8
$ws.Range("M74").Value = 66.85709999999995
$ws.Range("N74").Value = -3610.8
$ws.Range("H77").Value = 939.1
$ws.Range("I77").Value = 807.1429000000001
$ws.Range("J77").Value = 1862.8
$ws.Range("K77").Value = 4035.7145
$ws.Range("L77").Value = 9314
$ws.Range("M77").Value = 332.2855
$ws.Range("N77").Value = -18050
$ws.Range("H122").Value = 1528.8422
$ws.Range("I122").Value = 1591.0588
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 4773.1764
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -2323.1764
$ws.Range("N122").Value = -7900
$ws.Range("H136").Value = 1356.3959
$ws.Range("I136").Value = 1088.0322
$ws.Range("J136").Value = 1845.7646
$ws.Range("K136").Value = 3264.0966
$ws.Range("L136").Value = 5537.293799999999
$ws.Range("M136").Value = -714.0966000000003
$ws.Range("N136").Value = -10637.2938

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 60000
$ws.Range("J4").Value = 60000
$ws.Range("L4").Value = 60000
$ws.Range("N4").Value = -60224
$ws.Range("H14").Value = 1400
$ws.Range("J14").Value = 1400
$ws.Range("L14").Value = 1400
$ws.Range("N14").Value = -1740
$ws.Range("H58").Value = 1047.9231
$ws.Range("I58").Value = 887.6667
$ws.Range("J58").Value = 1408.5
$ws.Range("K58").Value = 887.6667
$ws.Range("L58").Value = 1408.5
$ws.Range("M58").Value = -684.6667
$ws.Range("N58").Value = -1814.5
$ws.Range("H94").Value = 66666960
$ws.Range("I94").Value = 250000220
$ws.Range("J94").Value = 318.0909
$ws.Range("K94").Value = 250000220
$ws.Range("L94").Value = 318.0909
$ws.Range("M94").Value = -249999769
$ws.Range("N94").Value = -1220.0909
$ws.Range("H132").Value = 2952.7058
$ws.Range("I132").Value = 3200.182
$ws.Range("J132").Value = 2499
$ws.Range("K132").Value = 9600.545999999998
$ws.Range("L132").Value = 7497
$ws.Range("M132").Value = -7070.545999999998
$ws.Range("N132").Value = -12557
$ws.Range("H136").Value = 1047.9231
$ws.Range("I136").Value = 887.6667
$ws.Range("J136").Value = 1408.5
$ws.Range("K136").Value = 2663.0001
$ws.Range("L136").Value = 4225.5
$ws.Range("M136").Value = -113.0001000000002
$ws.Range("N136").Value = -9325.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 100000000
$ws.Range("I4").Value = 100000000
$ws.Range("K4").Value = 300000000
$ws.Range("M4").Value = -299999888
$ws.Range("H5").Value = 1068
$ws.Range("J5").Value = 1213.875
$ws.Range("L5").Value = 3641.625
$ws.Range("N5").Value = -3865.625
$ws.Range("H10").Value = 22.714285
$ws.Range("I10").Value = 22.714285
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 68.142855
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 70.857145
$ws.Range("N10").ClearContents()
$ws.Range("H121").Value = 1077.7959
$ws.Range("I121").Value = 214.44444
$ws.Range("J121").Value = 1272.05
$ws.Range("K121").Value = 643.33332
$ws.Range("L121").Value = 3816.15
$ws.Range("M121").Value = 666.66668
$ws.Range("N121").Value = -6436.15
$ws.Range("H122").Value = 1383.5834
$ws.Range("I122").Value = 1714.8572
$ws.Range("J122").Value = 919.8
$ws.Range("K122").Value = 15433.7148
$ws.Range("L122").Value = 8278.199999999999
$ws.Range("M122").Value = -12983.7148
$ws.Range("N122").Value = -13178.2
$ws.Range("H131").Value = 2809.4038
$ws.Range("J131").Value = 955.3214
$ws.Range("L131").Value = 2865.9642
$ws.Range("N131").Value = -12945.9642
$ws.Range("H135").Value = 1068
$ws.Range("J135").Value = 1213.875
$ws.Range("L135").Value = 10924.875
$ws.Range("N135").Value = -15994.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 3350111
$ws.Range("I7").Value = 3350111
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 3350111
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -3349999
$ws.Range("N7").ClearContents()
$ws.Range("H8").Value = 3350111
$ws.Range("I8").Value = 3350111
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 3350111
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -3349972
$ws.Range("N8").ClearContents()
$ws.Range("H9").Value = 200
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()
$ws.Range("H12").Value = 5500017
$ws.Range("I12").Value = 6270020.5
$ws.Range("J12").Value = 1650000
$ws.Range("K12").Value = 6270020.5
$ws.Range("L12").Value = 1650000
$ws.Range("M12").Value = -6269880.5
$ws.Range("N12").Value = -1650280
$ws.Range("H109").Value = 10285
$ws.Range("J109").Value = 10285
$ws.Range("L109").Value = 10285
$ws.Range("N109").Value = -12365
$ws.Range("H113").Value = 41667680
$ws.Range("I113").Value = 250000000
$ws.Range("J113").Value = 1216
$ws.Range("K113").Value = 250000000
$ws.Range("L113").Value = 1216
$ws.Range("M113").Value = -249997830
$ws.Range("N113").Value = -5556

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 336311.12
$ws.Range("J2").Value = 8933.333000000001
$ws.Range("L2").Value = 8933.333000000001
$ws.Range("N2").Value = -9157.333000000001
$ws.Range("H22").Value = 1595.9166
$ws.Range("I22").Value = 1950.125
$ws.Range("J22").Value = 887.5
$ws.Range("K22").Value = 1950.125
$ws.Range("L22").Value = 887.5
$ws.Range("M22").Value = -1655.125
$ws.Range("N22").Value = -1477.5
$ws.Range("H27").Value = 1595.9166
$ws.Range("I27").Value = 1950.125
$ws.Range("J27").Value = 887.5
$ws.Range("K27").Value = 1950.125
$ws.Range("L27").Value = 887.5
$ws.Range("M27").Value = -1843.125
$ws.Range("N27").Value = -1101.5
$ws.Range("H40").Value = 1918.875
$ws.Range("I40").Value = 1300.1538
$ws.Range("J40").Value = 4600
$ws.Range("K40").Value = 1300.1538
$ws.Range("L40").Value = 4600
$ws.Range("M40").Value = -1164.1538
$ws.Range("N40").Value = -4872
$ws.Range("H122").Value = 2652.1035
$ws.Range("I122").Value = 2566.9524
$ws.Range("J122").Value = 2875.625
$ws.Range("K122").Value = 7700.8572
$ws.Range("L122").Value = 8626.875
$ws.Range("M122").Value = -5250.8572
$ws.Range("N122").Value = -13526.875
$ws.Range("H132").Value = 8410.034
$ws.Range("I132").Value = 19337.908
$ws.Range("J132").Value = 1731.8889
$ws.Range("K132").Value = 58013.724
$ws.Range("L132").Value = 5195.6667
$ws.Range("M132").Value = -55483.724
$ws.Range("N132").Value = -10255.6667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 450
$ws.Range("J2").Value = 450
$ws.Range("L2").Value = 450
$ws.Range("N2").Value = -674
$ws.Range("H13").Value = 500
$ws.Range("I13").Value = 500
$ws.Range("K13").Value = 500
$ws.Range("M13").Value = -360
$ws.Range("H132").Value = 1465.5333
$ws.Range("I132").Value = 1590.2858
$ws.Range("J132").Value = 1174.4445
$ws.Range("K132").Value = 4770.857400000001
$ws.Range("L132").Value = 3523.3335
$ws.Range("M132").Value = -2240.857400000001
$ws.Range("N132").Value = -8583.333500000001
$ws.Range("H136").Value = 3351.32
$ws.Range("I136").Value = 3859.5642
$ws.Range("J136").Value = 1549.3636
$ws.Range("K136").Value = 11578.6926
$ws.Range("L136").Value = 4648.0908
$ws.Range("M136").Value = -9028.692599999998
$ws.Range("N136").Value = -9748.0908
